# Rename the two repeated logo pictures (Pearson logo in the footers,
# BTec logo in the headers) that live in this document's header/footer
# parts. The edit swaps their default "imageN.ext" names:
#   PearsonLogo.png inline shapes: image1.png -> image2.png
#   BTec_Logo-Orange inline shapes: image2.jpg -> image1.jpg
# Every section's headers and footers (default / first-page / even-page,
# whichever exist) are walked so the rename is applied everywhere the
# logo appears, rather than hard-coding section/header indices.

$d = $word.ActiveDocument

function Rename-LogoShapes($storyRange) {
    $count = $storyRange.InlineShapes.Count
    for ($j = 1; $j -le $count; $j++) {
        $shape = $storyRange.InlineShapes.Item($j)
        $alt = $shape.AlternativeText

        if ($alt -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $shape.Name = "image2.png"
        } elseif ($alt -eq "BTec_Logo-Orange") {
            $shape.Name = "image1.jpg"
        }
    }
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $section = $d.Sections.Item($s)

    for ($i = 1; $i -le 3; $i++) {
        $hdr = $section.Headers.Item($i)
        if ($hdr.Exists) {
            Rename-LogoShapes($hdr.Range)
        }

        $ftr = $section.Footers.Item($i)
        if ($ftr.Exists) {
            Rename-LogoShapes($ftr.Range)
        }
    }
}
